# Updated cryptos list on Mon Sep 16 19:40:33 UTC 2024 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) figures for the cryptos sheet.
# Price values are digit/dot strings that Excel would otherwise coerce into numbers
# (dropping the thousands-separator dots or trailing zeros), so they are written with a
# leading apostrophe to force them to stay as text, matching the original inline-string data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.951.88"
$ws.Range("E2").Value = '  -3.22%  '
$ws.Range("D3").Value = "'2.288.09"
$ws.Range("E3").Value = '  -3.71%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'533.06"
$ws.Range("E5").Value = '  -4.28%  '
$ws.Range("D6").Value = "'130.96"
$ws.Range("E6").Value = '  -2.01%  '
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = "'0.584"
$ws.Range("E8").Value = '  -0.18%  '
$ws.Range("D9").Value = "'2.287.46"
$ws.Range("E9").Value = '  -3.71%  '
$ws.Range("D10").Value = "'0.0996"
$ws.Range("E10").Value = '  -5.91%  '
$ws.Range("E11").Value = '  -3.93%  '
$ws.Range("D12").Value = "'0.149"
$ws.Range("E12").Value = '  -0.36%  '
$ws.Range("E13").Value = '  -3.75%  '
$ws.Range("D14").Value = "'23.51"
$ws.Range("E14").Value = '  -3.87%  '
$ws.Range("D15").Value = "'2.698.49"
$ws.Range("D16").Value = "'57.895.17"
$ws.Range("E16").Value = '  -3.25%  '
$ws.Range("E17").Value = '  -4.80%  '
$ws.Range("D18").Value = "'2.303.06"
$ws.Range("E18").Value = '  -3.05%  '
$ws.Range("D19").Value = "'10.48"
$ws.Range("E19").Value = '  -5.75%  '
$ws.Range("E20").Value = '  -5.88%  '
$ws.Range("D21").Value = "'311.86"
$ws.Range("E21").Value = '  -2.97%  '
$ws.Range("D22").Value = "'6.38"
$ws.Range("E22").Value = '  -4.07%  '
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").Value = "'62.38"
$ws.Range("E24").Value = '  -2.73%  '
$ws.Range("E25").Value = '  -3.36%  '
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = '  -0.21%  '
$ws.Range("E27").Value = '  -5.49%  '
$ws.Range("E28").Value = '  -6.88%  '
$ws.Range("D29").Value = "'170.84"
$ws.Range("E30").Value = '  -6.28%  '
$ws.Range("E31").Value = '  -5.64%  '
$ws.Range("D32").Value = "'5.75"
$ws.Range("E32").Value = '  -5.42%  '
$ws.Range("E33").Value = '  -7.12%  '
$ws.Range("E34").Value = '  -5.58%  '
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("D36").Value = "'17.72"
$ws.Range("E36").Value = '  -2.34%  '
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("E38").Value = '  -7.28%  '
$ws.Range("D39").Value = "'3.89"
$ws.Range("D40").Value = "'38.15"
$ws.Range("E41").Value = '  -6.63%  '
$ws.Range("D42").Value = "'141.57"
$ws.Range("E42").Value = '  -2.16%  '
$ws.Range("D43").Value = "'289.17"
$ws.Range("E43").Value = '  -9.22%  '
$ws.Range("D44").Value = "'3.42"
$ws.Range("E44").Value = '  -3.23%  '
$ws.Range("D45").Value = "'0.0946"
$ws.Range("E45").Value = '  -2.28%  '
$ws.Range("D46").Value = "'0.0494"
$ws.Range("E46").Value = '  -3.34%  '
$ws.Range("E47").Value = '  -2.84%  '
$ws.Range("D48").Value = "'18.06"
$ws.Range("E48").Value = '  -8.17%  '
$ws.Range("E49").Value = '  -3.37%  '
$ws.Range("E51").Value = '  -0.69%  '
